# Applies edits described in commit "Make some progress on creating test data with the first 10 students"
# - Fills in survey answers for students in rows 5, 7-11 (C..J columns)
# - Updates C3 answer
# - Applies the "Good" (green) cell style to the corresponding schedule-availability columns
# - Moves the active selection to BC23

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update an existing answer cell, and fill in new answer cells for rows 5, 7, 8, 9, 10, 11 ---
$cellValues = @{
    "C3" = "A6"
    "C5" = "A8"
    "D5" = "A24"
    "H5" = "A1"
    "I5" = "A5"
    "C7" = "A2"
    "D7" = "A15"
    "E7" = "A18"
    "H7" = "A7"
    "I7" = "A21"
    "C8" = "A3"
    "D8" = "A10"
    "E8" = "A1"
    "F8" = "A19"
    "H8" = "A24"
    "I8" = "A15"
    "C9" = "A3"
    "D9" = "A20"
    "E9" = "A13"
    "H9" = "A23"
    "I9" = "A18"
    "J9" = "A10"
    "C10" = "A22"
    "D10" = "A11"
    "H10" = "A13"
    "I10" = "A12"
    "J10" = "A2"
    "C11" = "A11"
    "H11" = "A17"
    "I11" = "A22"
}

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}

# --- Apply the "Good" style to the highlighted availability cells ---
$styledCells = @("V5","W5","X5","Y5","AD5","AE5","AF5","AN5","AO5","AS5","AT5","AU5","AV5","AW5","BC5","BD5","BE5","BK5","BL5","M7","N7","O7","U7","V7","W7","AC7","AD7","AE7","AK7","AL7","AM7","AS7","AT7","AU7","BA7","BB7","BC7","BI7","BJ7","BK7","P8","Q8","R8","X8","Y8","Z8","AF8","AG8","AH8","AV8","AW8","AX8","BC8","BD8","BE8","BF8","BI8","BJ8","BK8","BL8","N9","O9","P9","V9","W9","X9","Y9","Z9","AD9","AE9","AF9","AG9","AH9","AL9","AM9","AN9","AO9","AP9","AT9","AU9","AV9","AW9","AX9","BB9","BC9","BD9","BE9","BF9","BJ9","BK9","BL9","O10","P10","Q10","X10","Y10","Z10","AE10","AF10","AG10","AM10","AN10","AO10","AU10","AV10","AW10","BC10","BD10","BE10","BK10","BL10","BM10","M11","N11","U11","V11","AC11","AD11","AK11","AL11","AS11","AT11","BA11","BB11","BC11","BD11","BE11","BF11","BI11","BJ11","BK11","BL11","BM11","BN11")

foreach ($addr in $styledCells) {
    $ws.Range($addr).Style = "Good"
}

# --- Move the active selection ---
$ws.Range("BC23").Select()
